$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header
$ws.Range("H1").Value = "Email"

# Update / fill in individual email addresses
$ws.Range("H2").Value  = "aidan.mccarron@mcaleer-rushe.co.uk"
$ws.Range("H3").Value  = "john.higgins@mcaleer-rushe.co.uk"
$ws.Range("H4").Value  = "michael.yohanis@mcaleer-rushe.co.uk"
$ws.Range("H5").Value  = "declan.mc@mcaleer-rushe.co.uk"
$ws.Range("H6").Value  = "lorcan.mulvey@mcaleer-rushe.co.uk"
$ws.Range("H7").Value  = "gerald.laverty@mcaleer-rushe.co.uk"
$ws.Range("H9").Value  = "steve.morris@mcaleer-rushe.co.uk"
$ws.Range("H10").Value = "lee.gray@mcaleer-rushe.co.uk"
$ws.Range("H11").Value = "eamonn.laverty@thorntonroofing.com"
$ws.Range("H12").Value = "eoin.gormley@mcaleer-rushe.co.uk"
$ws.Range("H13").Value = "paddy.connolly@mcaleer-rushe.co.uk"
$ws.Range("H14").Value = "daisy.butterworth@mcaleer-rushe.co.uk"
$ws.Range("H15").Value = "sinead.gorman@mcaleer-rushe.co.uk"
$ws.Range("H16").Value = "connor.graham@mcaleer-rushe.co.uk"
$ws.Range("H17").Value = "cathal.magee@mcaleer-rushe.co.uk"
$ws.Range("H19").Value = "nina.salandy@mcaleer-rushe.co.uk"
$ws.Range("H20").Value = "peter.coyle@mcaleer-rushe.co.uk"
$ws.Range("H21").Value = "orran.devine@mcaleer-rushe.co.uk"
$ws.Range("H22").Value = "niamh.heneghan@mcaleer-rushe.co.uk"

$wb.Save()
